# Update the Pakistan Reinsurance capital structure database values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell -> new value, applied identically to rows 2 and 3.
$values = @{
    "D" = 0.0542
    "E" = 0.0602
    "G" = 0.2128099173553719
    "H" = 0.2128099173553719
    "I" = 0.3347107438016529
    "J" = 0.2474767561983471
    "K" = 11.9
    "L" = 0.2458677685950413
    "M" = 3.35
    "N" = 0.06492248062015504
    "O" = 0.2815126050420168
    "P" = 3.35
    "Q" = 0.06492248062015504
    "R" = 0.2815126050420168
    "U" = 8.029999999999999
    "V" = 0.1556201550387597
    "W" = 0.1825153374233129
    "X" = 0.1475426634001101
    "Y" = 0.03497267402320281
    "Z" = 1.021097046413502
    "AA" = 0.2526977848101266
    "AB" = 0.1472519475124879
    "AC" = 0.1054458372976386
    "AD" = 0.154
    "AF" = 0.154
    "AG" = -7.875999999999999
    "AH" = 0.002975615411369169
    "AI" = 0.002526495389966204
    "AJ" = -0.1801299057725734
    "AK" = -0.1488171717935152
    "AL" = 0.008999999999999999
    "AM" = 0.008999999999999999
    "AN" = 0.009506172839506173
    "AO" = 1800
    "AP" = -0.4861728395061728
    "AQ" = 1800
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
